# Weekly update: insert a new data row at the top of the Zanahoria
# (Carrot) price records for "Agrícola del Norte S.A. de Arica", pushing
# all existing rows down by one, and fill it in with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 382; every row from 382..455 shifts
# down to 383..456 automatically.
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with this week's record.
$ws.Cells.Item(382, 1).Value = 1
$ws.Cells.Item(382, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(382, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(382, 4).Value = 44995
$ws.Cells.Item(382, 5).Value = 15
$ws.Cells.Item(382, 6).Value = 100114013
$ws.Cells.Item(382, 7).Value = "Zanahoria"
$ws.Cells.Item(382, 8).Value = "Sin especificar"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 110
$ws.Cells.Item(382, 11).Value = 9000
$ws.Cells.Item(382, 12).Value = 10000
$ws.Cells.Item(382, 13).Value = 9455
$ws.Cells.Item(382, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(382, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(382, 16).Value = 378
$ws.Cells.Item(382, 17).Value = 25
$ws.Cells.Item(382, 18).Value = "Hortaliza"
